$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'double[,]' 24,13
$arr[0,0] = 3.577770766035371
$arr[0,1] = 0.1606579634555771
$arr[0,2] = 0.457992860989819
$arr[0,3] = 0.1270481752441448
$arr[0,4] = 0
$arr[0,5] = 3.041145529798825
$arr[0,6] = 2.336121255631781
$arr[0,7] = 0
$arr[0,8] = 0.0377382600139029
$arr[0,9] = 0
$arr[0,10] = 0.6061498198585724
$arr[0,11] = 0
$arr[0,12] = 2.306366486997455
$arr[1,0] = 3.466044477585513
$arr[1,1] = 0.1416828904190197
$arr[1,2] = 0.4572801181569872
$arr[1,3] = 0.1275868731172256
$arr[1,4] = 0
$arr[1,5] = 3.02204000373294
$arr[1,6] = 2.335279542658583
$arr[1,7] = 0
$arr[1,8] = 0.03743876828796999
$arr[1,9] = 0
$arr[1,10] = 0.5989507204493805
$arr[1,11] = 0
$arr[1,12] = 2.330115751466614
$arr[2,0] = 3.399427577641234
$arr[2,1] = 0.1300697271083493
$arr[2,2] = 0.4570359368944139
$arr[2,3] = 0.1279539747464611
$arr[2,4] = 0
$arr[2,5] = 3.012067878838764
$arr[2,6] = 2.335891596710297
$arr[2,7] = 0
$arr[2,8] = 0.03725133294553906
$arr[2,9] = 0
$arr[2,10] = 0.5948308347769711
$arr[2,11] = 0
$arr[2,12] = 2.345460884204755
$arr[3,0] = 3.372778554349907
$arr[3,1] = 0.1253464934536623
$arr[3,2] = 0.4569850711601475
$arr[3,3] = 0.1281127207828234
$arr[3,4] = 0
$arr[3,5] = 3.008444895512213
$arr[3,6] = 2.33642436366145
$arr[3,7] = 0
$arr[3,8] = 0.03717405445250321
$arr[3,9] = 0
$arr[3,10] = 0.5932274275706817
$arr[3,11] = 0
$arr[3,12] = 2.35190571017948
$arr[4,0] = 3.368383552304408
$arr[4,1] = 0.1245627545663694
$arr[4,2] = 0.4569795626320712
$arr[4,3] = 0.1281396333662492
$arr[4,4] = 0
$arr[4,5] = 3.007869877584426
$arr[4,6] = 2.336529928045536
$arr[4,7] = 0
$arr[4,8] = 0.03716116808144321
$arr[4,9] = 0
$arr[4,10] = 0.5929657402754458
$arr[4,11] = 0
$arr[4,12] = 2.352987429683786
$arr[5,0] = 3.399066164733711
$arr[5,1] = 0.1300059908770663
$arr[5,2] = 0.4570350539618317
$arr[5,3] = 0.1279560785910707
$arr[5,4] = 0
$arr[5,5] = 3.012017235497922
$arr[5,6] = 2.335897635220164
$arr[5,7] = 0
$arr[5,8] = 0.03725029437993932
$arr[5,9] = 0
$arr[5,10] = 0.5948089051228322
$arr[5,11] = 0
$arr[5,12] = 2.345547026180746
$arr[6,0] = 3.538835294473245
$arr[6,1] = 0.1541074060316703
$arr[6,2] = 0.4577069537628944
$arr[6,3] = 0.1272263828728342
$arr[6,4] = 0
$arr[6,5] = 3.034191987683442
$arr[6,6] = 2.335596418376582
$arr[6,7] = 0
$arr[6,8] = 0.0376357279922459
$arr[6,9] = 0
$arr[6,10] = 0.6036051856904407
$arr[6,11] = 0
$arr[6,12] = 2.314396696199054
$arr[7,0] = 3.828719543756165
$arr[7,1] = 0.2016816178533247
$arr[7,2] = 0.4605602738313195
$arr[7,3] = 0.1260833141514528
$arr[7,4] = 0
$arr[7,5] = 3.09170643565011
$arr[7,6] = 2.343989099847249
$arr[7,7] = 0
$arr[7,8] = 0.03836369224300817
$arr[7,9] = 0
$arr[7,10] = 0.6232424887743235
$arr[7,11] = 0
$arr[7,12] = 2.25937756355065
$arr[8,0] = 4.05143968208688
$arr[8,1] = 0.2368458965893581
$arr[8,2] = 0.4635949376724255
$arr[8,3] = 0.1254184029025023
$arr[8,4] = 0
$arr[8,5] = 3.142628850442037
$arr[8,6] = 2.355672672609074
$arr[8,7] = 0
$arr[8,8] = 0.03888194628603792
$arr[8,9] = 0
$arr[8,10] = 0.6391343929521156
$arr[8,11] = 0
$arr[8,12] = 2.222669475966228
$arr[9,0] = 4.154902842127342
$arr[9,1] = 0.2528939206939356
$arr[9,2] = 0.4651797882350195
$arr[9,3] = 0.1251537790816304
$arr[9,4] = 0
$arr[9,5] = 3.167702112317727
$arr[9,6] = 2.36219524047803
$arr[9,7] = 0
$arr[9,8] = 0.03911420373928109
$arr[9,9] = 0
$arr[9,10] = 0.6466841258237821
$arr[9,11] = 0
$arr[9,12] = 2.206780131050785
$arr[10,0] = 4.194391980829835
$arr[10,1] = 0.2589786447439621
$arr[10,2] = 0.4658093503428518
$arr[10,3] = 0.1250590059852748
$arr[10,4] = 0
$arr[10,5] = 3.177473068821314
$arr[10,6] = 2.364839537096856
$arr[10,7] = 0
$arr[10,8] = 0.03920165696018785
$arr[10,9] = 0
$arr[10,10] = 0.6495892207802711
$arr[10,11] = 0
$arr[10,12] = 2.200879988359162
$arr[11,0] = 4.18587348637783
$arr[11,1] = 0.2576678444812615
$arr[11,2] = 0.4656724542639381
$arr[11,3] = 0.1250791755082119
$arr[11,4] = 0
$arr[11,5] = 3.175356402466036
$arr[11,6] = 2.364262275569558
$arr[11,7] = 0
$arr[11,8] = 0.03918284438841013
$arr[11,9] = 0
$arr[11,10] = 0.6489615013403096
$arr[11,11] = 0
$arr[11,12] = 2.202145487761292
$arr[12,0] = 4.158145419720256
$arr[12,1] = 0.253394359043881
$arr[12,2] = 0.4652309930404357
$arr[12,3] = 0.1251458731787682
$arr[12,4] = 0
$arr[12,5] = 3.168500427610866
$arr[12,6] = 2.362409290321096
$arr[12,7] = 0
$arr[12,8] = 0.03912140850984969
$arr[12,9] = 0
$arr[12,10] = 0.6469222037378586
$arr[12,11] = 0
$arr[12,12] = 2.206292379307492
$arr[13,0] = 4.141201570947544
$arr[13,1] = 0.2507777346215221
$arr[13,2] = 0.4649644166192957
$arr[13,3] = 0.125187434881278
$arr[13,4] = 0
$arr[13,5] = 3.164336976617449
$arr[13,6] = 2.36129700831151
$arr[13,7] = 0
$arr[13,8] = 0.03908371261788979
$arr[13,9] = 0
$arr[13,10] = 0.6456790914194954
$arr[13,11] = 0
$arr[13,12] = 2.208847694621355
$arr[14,0] = 4.044721418776192
$arr[14,1] = 0.2357981815048049
$arr[14,2] = 0.4634954785044698
$arr[14,3] = 0.1254364575720626
$arr[14,4] = 0
$arr[14,5] = 3.141028811537893
$arr[14,6] = 2.355270766529486
$arr[14,7] = 0
$arr[14,8] = 0.03886669777735641
$arr[14,9] = 0
$arr[14,10] = 0.6386474562089575
$arr[14,11] = 0
$arr[14,12] = 2.22372420423045
$arr[15,0] = 3.98608461713053
$arr[15,1] = 0.2266221166877074
$arr[15,2] = 0.4626466934770974
$arr[15,3] = 0.125598912666657
$arr[15,4] = 0
$arr[15,5] = 3.127220087700579
$arr[15,6] = 2.351883652093619
$arr[15,7] = 0
$arr[15,8] = 0.03873267394555846
$arr[15,9] = 0
$arr[15,10] = 0.6344159126700362
$arr[15,11] = 0
$arr[15,12] = 2.233058053455487
$arr[16,0] = 3.952560257996026
$arr[16,1] = 0.2213491292847038
$arr[16,2] = 0.4621777290601443
$arr[16,3] = 0.1256959154166086
$arr[16,4] = 0
$arr[16,5] = 3.119457191432673
$arr[16,6] = 2.350049101663927
$arr[16,7] = 0
$arr[16,8] = 0.03865525731554342
$arr[16,9] = 0
$arr[16,10] = 0.6320121945678778
$arr[16,11] = 0
$arr[16,12] = 2.238502790471067
$arr[17,0] = 3.941244140234687
$arr[17,1] = 0.2195646101154693
$arr[17,2] = 0.4620222485186929
$arr[17,3] = 0.1257293710637182
$arr[17,4] = 0
$arr[17,5] = 3.116859584931291
$arr[17,6] = 2.349447449584432
$arr[17,7] = 0
$arr[17,8] = 0.03862898858741204
$arr[17,9] = 0
$arr[17,10] = 0.6312035129633813
$arr[17,11] = 0
$arr[17,12] = 2.240359351566916
$arr[18,0] = 3.992305692255172
$arr[18,1] = 0.2275984213955837
$arr[18,2] = 0.4627350572381488
$arr[18,3] = 0.1255812503609448
$arr[18,4] = 0
$arr[18,5] = 3.128671457490384
$arr[18,6] = 2.352232451225547
$arr[18,7] = 0
$arr[18,8] = 0.03874697509073854
$arr[18,9] = 0
$arr[18,10] = 0.6348632462040911
$arr[18,11] = 0
$arr[18,12] = 2.232056564600221
$arr[19,0] = 4.16628140502354
$arr[19,1] = 0.2546493745973635
$arr[19,2] = 0.4653598624074817
$arr[19,3] = 0.1251261350379291
$arr[19,4] = 0
$arr[19,5] = 3.170506681361132
$arr[19,6] = 2.362948819902755
$arr[19,7] = 0
$arr[19,8] = 0.03913946717865535
$arr[19,9] = 0
$arr[19,10] = 0.6475199408934742
$arr[19,11] = 0
$arr[19,12] = 2.205071163194134
$arr[20,0] = 4.281791432651289
$arr[20,1] = 0.2723736601467976
$arr[20,2] = 0.4672467829490188
$arr[20,3] = 0.1248603619130453
$arr[20,4] = 0
$arr[20,5] = 3.199459522777914
$arr[20,6] = 2.370969111042939
$arr[20,7] = 0
$arr[20,8] = 0.03939308543782793
$arr[20,9] = 0
$arr[20,10] = 0.6560610149212351
$arr[20,11] = 0
$arr[20,12] = 2.188115689314849
$arr[21,0] = 4.219975821360151
$arr[21,1] = 0.2629096764016765
$arr[21,2] = 0.4662240003002154
$arr[21,3] = 0.1249993147827393
$arr[21,4] = 0
$arr[21,5] = 3.183858817710387
$arr[21,6] = 2.366595289038145
$arr[21,7] = 0
$arr[21,8] = 0.03925798793495083
$arr[21,9] = 0
$arr[21,10] = 0.6514778180436736
$arr[21,11] = 0
$arr[21,12] = 2.197102686632867
$arr[22,0] = 3.98949256269259
$arr[22,1] = 0.2271570264114757
$arr[22,2] = 0.4626950487613897
$arr[22,3] = 0.1255892242594996
$arr[22,4] = 0
$arr[22,5] = 3.128014745385769
$arr[22,6] = 2.352074408008662
$arr[22,7] = 0
$arr[22,8] = 0.0387405106787746
$arr[22,9] = 0
$arr[22,10] = 0.634660916228043
$arr[22,11] = 0
$arr[22,12] = 2.232509092976031
$arr[23,0] = 3.748595584739348
$arr[23,1] = 0.1887759749877773
$arr[23,2] = 0.4596237550901492
$arr[23,3] = 0.1263617885311756
$arr[23,4] = 0
$arr[23,5] = 3.074633621478853
$arr[23,6] = 2.340752743196816
$arr[23,7] = 0
$arr[23,8] = 0.03816970019543398
$arr[23,9] = 0
$arr[23,10] = 0.6176735645064895
$arr[23,11] = 0
$arr[23,12] = 2.273610593556917

$ws.Range("B2:N25").Value = $arr
